# Auto-update: GitHub Admin Log for delete-team
#
# Appends a new log row (row 9) to Sheet1, mirroring the existing
# "GitHub Admin Log" rows: Date & Time, Actions to perform, GitHub
# Organization, Team Name, Repository Name (if applicable), GitHub User
# Name, Permission Level, New Repo Name (if created new), Private Repo
# (True / False).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "2025-07-24 05:26:59"
$ws.Range("B9").Value = "delete-team"
$ws.Range("C9").Value = "new-organization97"
$ws.Range("D9").Value = "newteam"
$ws.Range("E9").Value = "demo"
$ws.Range("F9").Value = "GokulJ17"
# G9 (Permission Level) and H9 (New Repo Name) stay blank for this entry,
# same as e.g. row 3/4/7 above for actions that don't use them.

# Column I ("Private Repo (True / False)") stores the word "False" as
# literal TEXT in every existing row (I2:I8), not as a Boolean. Writing
# the bare word via .Value would make Excel auto-coerce it into a Boolean
# TRUE/FALSE value, so instead we build it as a text formula in a scratch
# cell and paste-special just the value into I9 - that keeps the result a
# plain text string "False", matching the rest of the column.
$ws.Range("K1").Formula = '="False"'
$ws.Range("K1").Copy()
$ws.Range("I9").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("K1").Value = ""
